# "restructure all folders complete"
# - Version sheet: bump Start Time / End Time to the next capture run
# - "3_3" / "1_8" sheets: replace each row's raw sample string and
#   recompute the Average/Max/Min/Sdev/Count columns from the new run

$wb = $excel.ActiveWorkbook

# --- Version sheet: Start Time / End Time ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = '23:09:56 03-02-2017 Thursday Pacific Standard Time'
$wsVersion.Range("C11").Value = '23:10:50 03-02-2017 Thursday Pacific Standard Time'

# --- "3_3" sheet: Deep Sleep / BT Idle / BT Page Scan rows ---
$ws33 = $wb.Worksheets.Item("3_3")

$ws33.Range("C3").Value = 0.011
$ws33.Range("D3").Value = 0.048
$ws33.Range("E3").Value = -0.011
$ws33.Range("F3").Value = 0.015
$ws33.Range("G3").Value = 98
$ws33.Range("H3").Value = '0.0,0.003628,0.028967,-0.004178,0.026087,0.003709,-0.002235,0.022004,0.005102,0.024226,-0.007982,0.029429,0.002595,0.022371,-0.000279,0.026746,-0.002697,0.025713,-0.00195,0.040381,0.00157,0.025441,0.009654,0.033879,0.002228,-0.001773,0.012989,-0.003716,0.025815,0.00409,0.033791,-0.006596,0.0064,0.029715,-0.00716,0.033981,0.00269,0.006406,0.03481,-0.001494,0.033417,0.002867,-9.5e-05,0.028587,0.007609,0.0,0.029056,-0.002975,0.031651,0.001108,0.03324,0.000374,0.024314,-0.003893,0.033519,-0.004742,0.029626,-0.000924,0.010014,0.006324,0.003166,-0.00214,-0.005,-0.010856,-0.001678,0.017643,0.03017,-0.000374,0.024511,0.007976,0.008825,-0.000469,0.001209,0.035727,0.007432,0.019402,0.003064,0.021345,0.000177,0.047901,-0.00659,0.029157,-0.001855,0.035557,0.009457,-0.008084,0.003804,-0.009572,0.027663,-0.005476,0.029898,0.005102,0.005285,0.003519,0.028967,-0.004273,0.017921,0.004178,0.034069,-0.003152'

$ws33.Range("C4").Value = 0.036
$ws33.Range("D4").Value = 0.062
$ws33.Range("E4").Value = -0.013
$ws33.Range("F4").Value = 0.014
$ws33.Range("G4").Value = 99
$ws33.Range("H4").Value = '0.0,-0.013091,-0.006128,0.007989,0.012357,0.011053,0.010679,0.011426,0.030829,0.026828,0.016712,0.009293,0.026359,0.040204,0.028872,0.02163,0.026087,0.036766,0.0404,0.033424,0.022194,0.039185,0.051617,0.05663,0.029443,0.032677,0.040666,0.050047,0.042248,0.028594,0.035469,0.046882,0.054877,0.039463,0.030544,0.040014,0.045591,0.050971,0.026644,0.046325,0.043166,0.047541,0.044198,0.026644,0.040197,0.048186,0.052833,0.031759,0.029905,0.036678,0.051624,0.040279,0.037146,0.030544,0.037133,0.061739,0.04159,0.030544,0.034443,0.050686,0.047534,0.045964,0.021726,0.043077,0.052453,0.042248,0.028791,0.036576,0.047629,0.052187,0.040109,0.029531,0.039457,0.054864,0.049661,0.024517,0.035007,0.043825,0.050598,0.036399,0.029993,0.030639,0.042615,0.043356,0.035652,0.03193,0.045027,0.051624,0.049851,0.029891,0.034633,0.04252,0.048471,0.040489,0.024144,0.030448,0.043268,0.049022,0.046141,0.029898'

$ws33.Range("C5").Value = -0.202
$ws33.Range("D5").Value = 0.047
$ws33.Range("E5").Value = -8.532999999999999
$ws33.Range("F5").Value = 1.265
$ws33.Range("G5").Value = 99
$ws33.Range("H5").Value = '0.0,0.032582,0.007602,0.036855,0.000924,0.015319,0.018193,0.03036,0.006868,0.042153,0.011046,0.001203,0.037507,0.038159,-8.532837,0.015781,0.008539,0.030645,0.008172,0.009654,0.040564,-0.004083,0.033791,0.015509,0.035279,0.002317,0.001671,0.034721,0.009015,0.001304,0.04661,0.008077,0.012629,0.031753,-0.000747,0.034348,0.00966,0.009565,0.035931,0.03788,0.005951,0.038342,0.003437,0.031019,0.00928,0.032215,0.005109,0.038995,0.00966,0.008635,0.036943,0.0118,-7.341299,0.033886,0.005115,0.004837,0.029429,0.013825,0.002595,0.012439,0.008077,0.033417,0.007609,0.036868,0.006678,0.002323,0.035177,0.0168,0.006501,0.032955,0.00269,0.007153,0.006216,0.007439,0.044742,0.014756,0.003532,0.002127,0.042602,0.006773,0.010958,0.005292,0.03676,0.017086,0.006318,0.019782,0.03964,0.002684,0.001481,0.039185,-6.119939,0.039552,0.00428,0.043261,0.005666,0.005,0.011793,0.041773,-0.001386,0.030822'

# --- "1_8" sheet: Deep Sleep / BT Idle / BT Page Scan rows ---
$ws18 = $wb.Worksheets.Item("1_8")

$ws18.Range("C3").Value = -0
$ws18.Range("D3").Value = 0.011
$ws18.Range("E3").Value = -0.008
$ws18.Range("F3").Value = 0.003
$ws18.Range("G3").Value = 99
$ws18.Range("H3").Value = '0.0,0.001986,0.003361,0.004717,-0.004144,-0.002177,-0.002732,-0.000247,-0.005266,-0.000506,-0.002183,-0.003848,-0.002466,-0.00024,-0.003552,-0.00053,0.00029,0.00394,0.00365,-0.000814,0.005309,-0.003058,0.000561,-0.00053,0.001696,0.003077,0.00558,4.9e-05,-0.003835,1.8e-05,0.001979,-0.001375,0.003095,-0.000795,-0.004162,0.002528,-0.000253,0.003342,-0.001634,-0.001091,0.003632,0.004495,-0.0055,-0.007739,0.000308,-0.003028,-0.001369,3.7e-05,-0.004964,-0.005235,-0.000802,0.005032,-0.001363,-0.002725,-0.001073,-0.004125,-0.001381,0.003897,0.003663,-0.001357,-0.005513,-0.004958,0.010853,0.003632,0.000327,-0.00304,-0.00246,0.004736,-0.000228,3.1e-05,-0.002183,-0.006351,-0.000216,0.000339,-0.002744,0.002799,0.001424,-0.000802,0.000857,0.001406,-0.004951,-0.000795,-0.000499,-0.001394,-0.001603,-0.003021,0.00251,0.00251,0.002528,0.000302,-0.004415,-0.001918,0.001424,-0.000518,0.001443,0.000302,0.002263,-0.001073,-0.001387,-0.001104'

$ws18.Range("C4").Value = 0.002
$ws18.Range("D4").Value = 0.011
$ws18.Range("E4").Value = -0.007
$ws18.Range("F4").Value = 0.003
$ws18.Range("G4").Value = 98
$ws18.Range("H4").Value = '0.0,0.005019,0.004452,0.006382,-0.001394,0.006394,0.0,0.004175,0.006388,0.001406,0.000561,-0.00135,0.001665,-0.000814,0.004736,0.007208,0.001653,0.007782,0.002232,0.004181,-0.002207,-0.002479,-0.001665,0.006962,0.006413,-0.004156,0.002485,-0.000561,0.004181,0.006382,-0.000845,-0.000277,0.000284,0.005556,0.004156,-0.001363,0.010569,0.002769,-0.001659,0.005284,0.000284,-0.000536,-0.001671,0.003058,0.002226,0.003885,0.000302,-0.003891,-0.00164,0.00751,0.000302,0.000271,0.001135,-0.001955,0.000839,-0.000555,0.0014,0.00193,0.00193,0.003077,-0.006672,0.007498,0.001949,0.000845,0.003897,-0.003065,0.003644,0.003601,0.004193,-0.000277,-0.001128,0.007239,0.00333,-0.001387,0.00891,0.004964,0.003077,0.001961,0.000845,0.000277,0.003866,0.001116,0.001646,-0.000555,0.00275,-0.00082,0.003872,0.000826,-0.003879,-0.001104,0.006413,-3.1e-05,0.002479,2.5e-05,0.000536,0.001973,-0.000808,0.006962,-0.000271,0.002522'

$ws18.Range("C5").Value = 0
$ws18.Range("D5").Value = 0.012
$ws18.Range("E5").Value = -0.008999999999999999
$ws18.Range("F5").Value = 0.003
$ws18.Range("G5").Value = 99
$ws18.Range("H5").Value = '0.0,0.000265,0.000567,-0.002263,0.000795,-0.000851,0.000253,0.00164,1.2e-05,-0.002547,-0.004199,-0.003903,0.006925,0.005556,0.000555,-0.001671,0.003872,-0.006142,0.000271,0.001924,-0.002812,0.006629,0.00275,-0.003354,-0.003058,0.005802,-0.005297,-0.002522,0.000265,-0.002781,-0.003626,0.005568,-0.004162,-0.000839,0.000247,0.000543,0.003866,0.000795,0.00304,-0.003922,0.000536,0.000259,-0.000592,0.003052,0.001375,-0.001979,-0.001683,-0.004193,0.003305,-0.004742,-0.000851,-0.002775,-0.003613,0.007467,-1.2e-05,0.003866,0.003576,0.000247,0.003589,-0.004723,-0.002503,0.000549,-0.003089,0.003305,0.002485,-0.003058,0.00193,-0.000857,0.000253,0.003046,0.000518,-0.000302,0.004452,0.011648,-0.00169,-1.2e-05,0.00053,0.004421,0.001085,-0.004452,0.002473,1.2e-05,0.000253,0.001659,0.003311,0.001665,0.005241,-0.005858,0.000814,0.003028,0.001091,-0.00029,0.000555,-0.008892,-0.001967,0.00497,-1.2e-05,-0.001135,0.001375,-0.000561'
